$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.738.17'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +1.36%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.727.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.27%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''0.9979'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  -0.14%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''240.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -0.76%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.9984'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -0.13%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.4823'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -0.97%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.2585'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -0.18%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.06182'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -0.03%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''1.726.88'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +0.26%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''15.86'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +2.36%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.06865'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -1.59%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''0.6037'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +1.14%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''4.463'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -1.20%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''76.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -0.19%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''0.9987'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -0.07%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''26.559.98'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = '''0.9983'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -0.15%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''0.000007146'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -0.53%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''11.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +0.50%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''1.946.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +0.05%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''4.418'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -0.29%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''8.537'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +0.49%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''5.053'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -0.68%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''139.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +1.21%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''15.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -0.16%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''1.776'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +3.09%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''106.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +0.19%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''1.371'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -2.11%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''4.008'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +2.74%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''0.07933'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -0.89%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''3.663'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +0.44%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  +0.48%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''2.597'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = '''0.9994'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +0.36%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''0.6175'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -0.72%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''0.9336'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +0.34%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''2.451'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +2.83%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''1.995'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +1.88%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.9978'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -0.10%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.01497'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +1.49%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''5.608'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +3.14%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''99.72'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -0.79%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.3829'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -0.14%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''6.785'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -1.30%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.1154'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -0.78%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.05359'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = '''7.892'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +2.70%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''30.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -0.32%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''1.241'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +1.84%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''51.41'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +0.86%  '
$ws.Range("E51").Style = "Normal"
